$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Dades_Període")
$ws2 = $wb.Worksheets.Item("Estudi_Capçaleres")

# --- Column H (DATA_EXTRACCIO) updates for all data rows ---
$hValues = @{
    2 = "2026-02-20 06:31:32"
    3 = "2026-02-20 06:31:33"
    4 = "2026-02-20 06:31:33"
    5 = "2026-02-20 06:31:33"
    6 = "2026-02-20 06:31:33"
    7 = "2026-02-20 06:31:36"
    8 = "2026-02-20 06:31:37"
    9 = "2026-02-20 06:31:37"
    10 = "2026-02-20 06:31:37"
    11 = "2026-02-20 06:31:37"
    12 = "2026-02-20 06:31:40"
    13 = "2026-02-20 06:31:41"
    14 = "2026-02-20 06:31:41"
    15 = "2026-02-20 06:31:41"
    16 = "2026-02-20 06:31:41"
    17 = "2026-02-20 06:31:44"
    18 = "2026-02-20 06:31:45"
    19 = "2026-02-20 06:31:45"
    20 = "2026-02-20 06:31:45"
    21 = "2026-02-20 06:31:45"
    22 = "2026-02-20 06:31:48"
    23 = "2026-02-20 06:31:49"
    24 = "2026-02-20 06:31:49"
    25 = "2026-02-20 06:31:49"
    26 = "2026-02-20 06:31:49"
    27 = "2026-02-20 06:31:52"
    28 = "2026-02-20 06:31:53"
    29 = "2026-02-20 06:31:53"
    30 = "2026-02-20 06:31:53"
    31 = "2026-02-20 06:31:53"
    32 = "2026-02-20 06:31:56"
    33 = "2026-02-20 06:31:57"
    34 = "2026-02-20 06:31:57"
    35 = "2026-02-20 06:31:57"
    36 = "2026-02-20 06:31:57"
    37 = "2026-02-20 06:32:00"
    38 = "2026-02-20 06:32:01"
    39 = "2026-02-20 06:32:01"
    40 = "2026-02-20 06:32:01"
    41 = "2026-02-20 06:32:01"
    42 = "2026-02-20 06:32:03"
    43 = "2026-02-20 06:32:05"
    44 = "2026-02-20 06:32:05"
    45 = "2026-02-20 06:32:05"
    46 = "2026-02-20 06:32:05"
    47 = "2026-02-20 06:32:07"
    48 = "2026-02-20 06:32:09"
    49 = "2026-02-20 06:32:09"
    50 = "2026-02-20 06:32:09"
    51 = "2026-02-20 06:32:09"
    52 = "2026-02-20 06:32:11"
    53 = "2026-02-20 06:32:13"
    54 = "2026-02-20 06:32:13"
    55 = "2026-02-20 06:32:13"
    56 = "2026-02-20 06:32:13"
    57 = "2026-02-20 06:32:15"
    58 = "2026-02-20 06:32:17"
    59 = "2026-02-20 06:32:17"
    60 = "2026-02-20 06:32:17"
    61 = "2026-02-20 06:32:17"
    62 = "2026-02-20 06:32:19"
    63 = "2026-02-20 06:32:20"
    64 = "2026-02-20 06:32:20"
    65 = "2026-02-20 06:32:20"
    66 = "2026-02-20 06:32:20"
    67 = "2026-02-20 06:32:23"
    68 = "2026-02-20 06:32:24"
    69 = "2026-02-20 06:32:24"
    70 = "2026-02-20 06:32:24"
    71 = "2026-02-20 06:32:24"
    72 = "2026-02-20 06:32:26"
    73 = "2026-02-20 06:32:28"
    74 = "2026-02-20 06:32:28"
    75 = "2026-02-20 06:32:28"
    76 = "2026-02-20 06:32:28"
    77 = "2026-02-20 06:32:30"
    78 = "2026-02-20 06:32:32"
    79 = "2026-02-20 06:32:32"
    80 = "2026-02-20 06:32:32"
    81 = "2026-02-20 06:32:32"
    82 = "2026-02-20 06:32:35"
    83 = "2026-02-20 06:32:36"
    84 = "2026-02-20 06:32:36"
    85 = "2026-02-20 06:32:36"
    86 = "2026-02-20 06:32:36"
    87 = "2026-02-20 06:32:39"
    88 = "2026-02-20 06:32:40"
    89 = "2026-02-20 06:32:40"
    90 = "2026-02-20 06:32:40"
    91 = "2026-02-20 06:32:40"
    92 = "2026-02-20 06:32:42"
    93 = "2026-02-20 06:32:43"
    94 = "2026-02-20 06:32:43"
    95 = "2026-02-20 06:32:43"
    96 = "2026-02-20 06:32:43"
    97 = "2026-02-20 06:32:46"
    98 = "2026-02-20 06:32:47"
    99 = "2026-02-20 06:32:47"
    100 = "2026-02-20 06:32:47"
    101 = "2026-02-20 06:32:47"
    102 = "2026-02-20 06:32:50"
    103 = "2026-02-20 06:32:51"
    104 = "2026-02-20 06:32:51"
    105 = "2026-02-20 06:32:51"
    106 = "2026-02-20 06:32:51"
    107 = "2026-02-20 06:32:54"
    108 = "2026-02-20 06:32:55"
    109 = "2026-02-20 06:32:55"
    110 = "2026-02-20 06:32:55"
    111 = "2026-02-20 06:32:55"
    112 = "2026-02-20 06:32:58"
    113 = "2026-02-20 06:32:59"
    114 = "2026-02-20 06:32:59"
    115 = "2026-02-20 06:32:59"
    116 = "2026-02-20 06:32:59"
    117 = "2026-02-20 06:33:01"
    118 = "2026-02-20 06:33:03"
    119 = "2026-02-20 06:33:03"
    120 = "2026-02-20 06:33:03"
    121 = "2026-02-20 06:33:03"
    122 = "2026-02-20 06:33:05"
    123 = "2026-02-20 06:33:07"
    124 = "2026-02-20 06:33:07"
    125 = "2026-02-20 06:33:07"
    126 = "2026-02-20 06:33:07"
    127 = "2026-02-20 06:33:09"
    128 = "2026-02-20 06:33:10"
    129 = "2026-02-20 06:33:10"
    130 = "2026-02-20 06:33:10"
    131 = "2026-02-20 06:33:10"
    132 = "2026-02-20 06:33:13"
    133 = "2026-02-20 06:33:14"
    134 = "2026-02-20 06:33:14"
    135 = "2026-02-20 06:33:14"
    136 = "2026-02-20 06:33:14"
    137 = "2026-02-20 06:33:17"
    138 = "2026-02-20 06:33:18"
    139 = "2026-02-20 06:33:18"
    140 = "2026-02-20 06:33:18"
    141 = "2026-02-20 06:33:18"
    142 = "2026-02-20 06:33:21"
    143 = "2026-02-20 06:33:22"
    144 = "2026-02-20 06:33:22"
    145 = "2026-02-20 06:33:22"
    146 = "2026-02-20 06:33:22"
    147 = "2026-02-20 06:33:25"
    148 = "2026-02-20 06:33:26"
    149 = "2026-02-20 06:33:26"
    150 = "2026-02-20 06:33:26"
    151 = "2026-02-20 06:33:26"
    152 = "2026-02-20 06:33:29"
    153 = "2026-02-20 06:33:30"
    154 = "2026-02-20 06:33:30"
    155 = "2026-02-20 06:33:30"
    156 = "2026-02-20 06:33:30"
    157 = "2026-02-20 06:33:32"
    158 = "2026-02-20 06:33:34"
    159 = "2026-02-20 06:33:34"
    160 = "2026-02-20 06:33:34"
    161 = "2026-02-20 06:33:34"
    162 = "2026-02-20 06:33:36"
    163 = "2026-02-20 06:33:38"
    164 = "2026-02-20 06:33:38"
    165 = "2026-02-20 06:33:38"
    166 = "2026-02-20 06:33:38"
    167 = "2026-02-20 06:33:40"
    168 = "2026-02-20 06:33:42"
    169 = "2026-02-20 06:33:42"
    170 = "2026-02-20 06:33:42"
    171 = "2026-02-20 06:33:42"
    172 = "2026-02-20 06:33:44"
    173 = "2026-02-20 06:33:45"
    174 = "2026-02-20 06:33:45"
    175 = "2026-02-20 06:33:45"
    176 = "2026-02-20 06:33:45"
    177 = "2026-02-20 06:33:48"
    178 = "2026-02-20 06:33:49"
    179 = "2026-02-20 06:33:49"
    180 = "2026-02-20 06:33:49"
    181 = "2026-02-20 06:33:49"
    182 = "2026-02-20 06:33:52"
    183 = "2026-02-20 06:33:53"
    184 = "2026-02-20 06:33:53"
    185 = "2026-02-20 06:33:53"
    186 = "2026-02-20 06:33:53"
    187 = "2026-02-20 06:33:56"
    188 = "2026-02-20 06:33:57"
    189 = "2026-02-20 06:33:57"
    190 = "2026-02-20 06:33:57"
    191 = "2026-02-20 06:33:57"
    192 = "2026-02-20 06:33:59"
    193 = "2026-02-20 06:34:01"
    194 = "2026-02-20 06:34:01"
    195 = "2026-02-20 06:34:01"
    196 = "2026-02-20 06:34:01"
    197 = "2026-02-20 06:34:03"
    198 = "2026-02-20 06:34:05"
    199 = "2026-02-20 06:34:05"
    200 = "2026-02-20 06:34:05"
    201 = "2026-02-20 06:34:05"
    202 = "2026-02-20 06:34:07"
    203 = "2026-02-20 06:34:09"
    204 = "2026-02-20 06:34:09"
    205 = "2026-02-20 06:34:09"
    206 = "2026-02-20 06:34:09"
    207 = "2026-02-20 06:34:11"
    208 = "2026-02-20 06:34:13"
    209 = "2026-02-20 06:34:13"
    210 = "2026-02-20 06:34:13"
    211 = "2026-02-20 06:34:13"
    212 = "2026-02-20 06:34:15"
    213 = "2026-02-20 06:34:17"
    214 = "2026-02-20 06:34:17"
    215 = "2026-02-20 06:34:17"
    216 = "2026-02-20 06:34:17"
    217 = "2026-02-20 06:34:19"
    218 = "2026-02-20 06:34:21"
    219 = "2026-02-20 06:34:21"
    220 = "2026-02-20 06:34:21"
    221 = "2026-02-20 06:34:21"
    222 = "2026-02-20 06:34:23"
    223 = "2026-02-20 06:34:24"
    224 = "2026-02-20 06:34:24"
    225 = "2026-02-20 06:34:24"
    226 = "2026-02-20 06:34:24"
}
foreach ($r in $hValues.Keys) {
    $ws1.Cells.Item($r, 8).Value = $hValues[$r]
}

# --- Column I (HORA_CONSULTA_UTC) and J (URL_FONT) updates on first row of each station block ---
$ijValues = @{
    2 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YT&dia=2026-02-20T06:00Z")
    7 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z1&dia=2026-02-20T06:00Z")
    12 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=DN&dia=2026-02-20T06:00Z")
    17 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z6&dia=2026-02-20T06:00Z")
    22 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=DJ&dia=2026-02-20T06:00Z")
    27 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=X4&dia=2026-02-20T06:00Z")
    32 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=D5&dia=2026-02-20T06:00Z")
    37 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YS&dia=2026-02-20T06:00Z")
    42 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=UN&dia=2026-02-20T06:00Z")
    47 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=MS&dia=2026-02-20T06:00Z")
    52 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=W1&dia=2026-02-20T06:00Z")
    57 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=DP&dia=2026-02-20T06:00Z")
    62 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XL&dia=2026-02-20T06:00Z")
    67 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=VZ&dia=2026-02-20T06:00Z")
    72 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z7&dia=2026-02-20T06:00Z")
    77 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XK&dia=2026-02-20T06:00Z")
    82 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T06:00Z")
    87 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YU&dia=2026-02-20T06:00Z")
    92 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=ZE&dia=2026-02-20T06:00Z")
    97 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=CD&dia=2026-02-20T06:00Z")
    102 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z2&dia=2026-02-20T06:00Z")
    107 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z5&dia=2026-02-20T06:00Z")
    112 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=VK&dia=2026-02-20T06:00Z")
    117 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z3&dia=2026-02-20T06:00Z")
    122 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=CG&dia=2026-02-20T06:00Z")
    127 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=Z9&dia=2026-02-20T06:00Z")
    132 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YB&dia=2026-02-20T06:00Z")
    137 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YP&dia=2026-02-20T06:00Z")
    142 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=J5&dia=2026-02-20T06:00Z")
    147 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=D6&dia=2026-02-20T06:00Z")
    152 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XR&dia=2026-02-20T06:00Z")
    157 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YA&dia=2026-02-20T06:00Z")
    162 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=DG&dia=2026-02-20T06:00Z")
    167 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=X5&dia=2026-02-20T06:00Z")
    172 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=D4&dia=2026-02-20T06:00Z")
    177 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=CI&dia=2026-02-20T06:00Z")
    182 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XS&dia=2026-02-20T06:00Z")
    187 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=ZC&dia=2026-02-20T06:00Z")
    192 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XH&dia=2026-02-20T06:00Z")
    197 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XE&dia=2026-02-20T06:00Z")
    202 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=UE&dia=2026-02-20T06:00Z")
    207 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=XO&dia=2026-02-20T06:00Z")
    212 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=VS&dia=2026-02-20T06:00Z")
    217 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=YN&dia=2026-02-20T06:00Z")
    222 = @("06:00", "https://www.meteo.cat/observacions/xema/dades?codi=D7&dia=2026-02-20T06:00Z")
}
foreach ($r in $ijValues.Keys) {
    $ws1.Cells.Item($r, 9).Value = $ijValues[$r][0]
    $ws1.Cells.Item($r, 10).Value = $ijValues[$r][1]
}

# --- Columns S, AJ, AR (precipitation) corrections to 0.0 on specific rows ---
$parZeroRows = @(8, 9, 18, 19, 20, 21, 74, 106, 109, 118, 119, 188, 189, 190, 213, 214)
foreach ($r in $parZeroRows) {
    $ws1.Cells.Item($r, 19).Value = "0.0"   # S
    $ws1.Cells.Item($r, 36).Value = "0.0"   # AJ
    $ws1.Cells.Item($r, 44).Value = "0.0"   # AR
}

# --- Sheet2 (Estudi_Capçaleres) column F (URL_FONT) updates ---
$fValues = @{
    2 = "https://www.meteo.cat/observacions/xema/dades?codi=YT&dia=2026-02-20T06:00Z"
    3 = "https://www.meteo.cat/observacions/xema/dades?codi=Z1&dia=2026-02-20T06:00Z"
    4 = "https://www.meteo.cat/observacions/xema/dades?codi=DN&dia=2026-02-20T06:00Z"
    5 = "https://www.meteo.cat/observacions/xema/dades?codi=Z6&dia=2026-02-20T06:00Z"
    6 = "https://www.meteo.cat/observacions/xema/dades?codi=DJ&dia=2026-02-20T06:00Z"
    7 = "https://www.meteo.cat/observacions/xema/dades?codi=X4&dia=2026-02-20T06:00Z"
    8 = "https://www.meteo.cat/observacions/xema/dades?codi=D5&dia=2026-02-20T06:00Z"
    9 = "https://www.meteo.cat/observacions/xema/dades?codi=YS&dia=2026-02-20T06:00Z"
    10 = "https://www.meteo.cat/observacions/xema/dades?codi=UN&dia=2026-02-20T06:00Z"
    11 = "https://www.meteo.cat/observacions/xema/dades?codi=MS&dia=2026-02-20T06:00Z"
    12 = "https://www.meteo.cat/observacions/xema/dades?codi=W1&dia=2026-02-20T06:00Z"
    13 = "https://www.meteo.cat/observacions/xema/dades?codi=DP&dia=2026-02-20T06:00Z"
    14 = "https://www.meteo.cat/observacions/xema/dades?codi=XL&dia=2026-02-20T06:00Z"
    15 = "https://www.meteo.cat/observacions/xema/dades?codi=VZ&dia=2026-02-20T06:00Z"
    16 = "https://www.meteo.cat/observacions/xema/dades?codi=Z7&dia=2026-02-20T06:00Z"
    17 = "https://www.meteo.cat/observacions/xema/dades?codi=XK&dia=2026-02-20T06:00Z"
    18 = "https://www.meteo.cat/observacions/xema/dades?codi=XJ&dia=2026-02-20T06:00Z"
    19 = "https://www.meteo.cat/observacions/xema/dades?codi=YU&dia=2026-02-20T06:00Z"
    20 = "https://www.meteo.cat/observacions/xema/dades?codi=ZE&dia=2026-02-20T06:00Z"
    21 = "https://www.meteo.cat/observacions/xema/dades?codi=CD&dia=2026-02-20T06:00Z"
    22 = "https://www.meteo.cat/observacions/xema/dades?codi=Z2&dia=2026-02-20T06:00Z"
    23 = "https://www.meteo.cat/observacions/xema/dades?codi=Z5&dia=2026-02-20T06:00Z"
    24 = "https://www.meteo.cat/observacions/xema/dades?codi=VK&dia=2026-02-20T06:00Z"
    25 = "https://www.meteo.cat/observacions/xema/dades?codi=Z3&dia=2026-02-20T06:00Z"
    26 = "https://www.meteo.cat/observacions/xema/dades?codi=CG&dia=2026-02-20T06:00Z"
    27 = "https://www.meteo.cat/observacions/xema/dades?codi=Z9&dia=2026-02-20T06:00Z"
    28 = "https://www.meteo.cat/observacions/xema/dades?codi=YB&dia=2026-02-20T06:00Z"
    29 = "https://www.meteo.cat/observacions/xema/dades?codi=YP&dia=2026-02-20T06:00Z"
    30 = "https://www.meteo.cat/observacions/xema/dades?codi=J5&dia=2026-02-20T06:00Z"
    31 = "https://www.meteo.cat/observacions/xema/dades?codi=D6&dia=2026-02-20T06:00Z"
    32 = "https://www.meteo.cat/observacions/xema/dades?codi=XR&dia=2026-02-20T06:00Z"
    33 = "https://www.meteo.cat/observacions/xema/dades?codi=YA&dia=2026-02-20T06:00Z"
    34 = "https://www.meteo.cat/observacions/xema/dades?codi=DG&dia=2026-02-20T06:00Z"
    35 = "https://www.meteo.cat/observacions/xema/dades?codi=X5&dia=2026-02-20T06:00Z"
    36 = "https://www.meteo.cat/observacions/xema/dades?codi=D4&dia=2026-02-20T06:00Z"
    37 = "https://www.meteo.cat/observacions/xema/dades?codi=CI&dia=2026-02-20T06:00Z"
    38 = "https://www.meteo.cat/observacions/xema/dades?codi=XS&dia=2026-02-20T06:00Z"
    39 = "https://www.meteo.cat/observacions/xema/dades?codi=ZC&dia=2026-02-20T06:00Z"
    40 = "https://www.meteo.cat/observacions/xema/dades?codi=XH&dia=2026-02-20T06:00Z"
    41 = "https://www.meteo.cat/observacions/xema/dades?codi=XE&dia=2026-02-20T06:00Z"
    42 = "https://www.meteo.cat/observacions/xema/dades?codi=UE&dia=2026-02-20T06:00Z"
    43 = "https://www.meteo.cat/observacions/xema/dades?codi=XO&dia=2026-02-20T06:00Z"
    44 = "https://www.meteo.cat/observacions/xema/dades?codi=VS&dia=2026-02-20T06:00Z"
    45 = "https://www.meteo.cat/observacions/xema/dades?codi=YN&dia=2026-02-20T06:00Z"
    46 = "https://www.meteo.cat/observacions/xema/dades?codi=D7&dia=2026-02-20T06:00Z"
}
foreach ($r in $fValues.Keys) {
    $ws2.Cells.Item($r, 6).Value = $fValues[$r]
}